$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added EXTERNAL_IDENTIFIERS (ID1 / ID2) columns with sample values
$ws.Range("F1").Value = "ID1"
$ws.Range("F2").Value = "id1-1"
$ws.Range("F3").Value = "id1-2"
$ws.Range("F4").Value = "id1-3"

$ws.Range("G1").Value = "ID2"
$ws.Range("G2").Value = "id2-1"
$ws.Range("G3").Value = "id2-2"
$ws.Range("G4").Value = "id2-3"

# Move the selection onto the newly populated column
$ws.Range("F1").Select()
